# "hand counts added to egg notes"
# Fill in the B-column (hand count) values that were recorded for this
# batch of egg-count notes, and move the sheet's active selection to C10
# (clearing the previous scroll-locked topLeftCell/selection at B46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 previously held the "NA" placeholder text; it now has a real hand count.
$ws.Range("B8").Value = 224

# Newly hand-counted rows (B29:B46) that had no count before.
$ws.Range("B29").Value = 715
$ws.Range("B30").Value = 120
$ws.Range("B31").Value = 830
$ws.Range("B33").Value = 157
$ws.Range("B35").Value = 853
$ws.Range("B36").Value = 526
$ws.Range("B37").Value = 195
$ws.Range("B38").Value = 368
$ws.Range("B39").Value = 338
$ws.Range("B40").Value = 85
$ws.Range("B41").Value = 647
$ws.Range("B43").Value = 841
$ws.Range("B44").Value = 770
$ws.Range("B45").Value = 747
$ws.Range("B46").Value = 326

# Move the selection/viewport off of B46 and onto C10 (also drops the
# stale topLeftCell="A32" scroll anchor).
$ws.Range("C10").Select()
